# Generate Report for handback
# Updates the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps on the last data row
# of the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-22 03:02:39"
$wsZhCn.Range("G5").Value = "2016-01-22 03:03:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-22 03:02:53"
$wsDeDe.Range("G5").Value = "2016-01-22 03:03:49"
